# Actualización automática del tracker
# Rellena las columnas "resultado" (G) y "profit" (H) para las filas
# que ya tienen un pronóstico resuelto (filas 3 y 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Fallo"
$ws.Range("H3").Value = -1

$ws.Range("G5").Value = "Acierto"
$ws.Range("H5").Value = 1.75
